$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the typo "hace" -> "have" and reorder/extend the "Options" intent sentences.
# Row 112 currently holds "What options do I hace?" -> replace with "What options do I have?"
$ws.Range("B112").Value = "What options do I have?"
# Row 113 currently holds "Tell me what can I ask for" -> keep same text, now sourced from new position
$ws.Range("B113").Value = "Tell me what can I ask for"
# Row 114 currently holds "Show me what you got" -> keep same text, now sourced from new position
$ws.Range("B114").Value = "Show me what you got"

# New rows for the "Options" intent
$ws.Range("A115").Value = "Options"
$ws.Range("B115").Value = "How many options do I have?"

$ws.Range("A116").Value = "Options"
$ws.Range("B116").Value = "Is there anything more I can do?"

# New rows for the "Headers" intent
$ws.Range("A117").Value = "Headers"
$ws.Range("B117").Value = "Tell me everything about Data Science"

$ws.Range("A118").Value = "Headers"
$ws.Range("B118").Value = "Show me more of Barcelona"

$ws.Range("A119").Value = "Headers"
$ws.Range("B119").Value = "What headers does this article have?"

$ws.Range("A120").Value = "Headers"
$ws.Range("B120").Value = "Give me more information please"

$ws.Range("A121").Value = "Headers"
$ws.Range("B121").Value = "I want to know more of Napoleon"

$ws.Range("B122").Select()
